$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Column I (No_Hp) switches from text-with-leading-zero to a plain
#    number for every existing row, and row 6 gets a new phone number.
#    A cell that already holds a shared-string has to be Clear()'d
#    before it will persist as a genuine <v> number instead of being
#    re-interned as a string, so we clear first, assign the numbers,
#    then restore the original "s=1" cell style (copied from a sibling
#    cell that already carries that style) without touching its value.
# ---------------------------------------------------------------------
foreach ($r in 1..6) {
    $ws.Range("I$r").Clear()
}
$ws.Range("I1").Value = 89237812378
$ws.Range("I2").Value = 89237812377
$ws.Range("I3").Value = 89237812379
$ws.Range("I4").Value = 89237812380
$ws.Range("I5").Value = 89237812381
$ws.Range("I6").Value = 89237812382

$ws.Range("E1").Copy()
$ws.Range("I1:I6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Row 6's e-mail column gets corrected to Ade Yuliana's mailbox
#    (the faculty row underneath points at Dini's old address instead,
#    mirroring the same display/value mismatch that already exists
#    elsewhere in this sheet).
# ---------------------------------------------------------------------
$ws.Range("H6").Value = "adeYul@poltekedc.ac.id"

# ---------------------------------------------------------------------
# 3) New row 7: Ade Yuliana, M.T.
#    Columns that must stay textual (leading zeros / huge ID strings /
#    literal date text) are pre-formatted as Text so Excel doesn't
#    "helpfully" coerce them into numbers or dates.
# ---------------------------------------------------------------------
$ws.Range("A7").NumberFormat = "@"
$ws.Range("B7").NumberFormat = "@"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("F7").NumberFormat = "@"

$ws.Range("A7").Value = "0003107907"
$ws.Range("B7").Value = "3217023989012396"
$ws.Range("C7").Value = "197910032005012007"
$ws.Range("D7").Value = "Ade Yuliana, M.T."
$ws.Range("E7").Value = "Bandung"
$ws.Range("F7").Value = "2000-10-06"
$ws.Range("G7").Value = "Islam"
$ws.Range("H7").Value = "diniR@poltekedc.ac.id"
$ws.Range("I7").Value = 89237812383
$ws.Range("J7").Value = "Jl. Bandung"
$ws.Range("K7").Value = 10
$ws.Range("L7").Value = "Aktif"
$ws.Range("M7").Value = "Perempuan"

# Hyperlink for the new row's e-mail cell. Add it with the display
# text first (so the relationship's cached display text is right),
# then put the real per-row address back into the cell - exactly the
# same display/value split already used for H3:H6 above.
$ws.Hyperlinks.Add($ws.Range("H7"), "mailto:aris@poltekedc.ac.id", "", "", "aris@poltekedc.ac.id")
$ws.Range("H7").Value = "diniR@poltekedc.ac.id"

# Match row 7's look (borders/fonts/number formats) to row 6, the
# previous last row of the table - done after the hyperlink so Excel's
# automatic hyperlink styling on H7 gets replaced by the real column
# style again.
$ws.Range("A6:M6").Copy()
$ws.Range("A7:M7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4) Recompute the autofit column widths now that the data has changed,
#    then move the selection the same way the author left it.
# ---------------------------------------------------------------------
$ws.Columns.AutoFit()
$ws.Range("N7").Select()
